{"js": "// Replace the first code listing (\"has_unique_chars\") with the new\n// \"bank\"/\"client\" class listing, and drop the two blank paragraphs plus\n// the duplicate \"has_unique_chars\" listing that followed it, leaving the\n// trailing blank paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: first \"has_unique_chars\" listing -> becomes the new bank code.\n// Paragraph 1: blank paragraph -> removed.\n// Paragraph 2: blank paragraph -> removed.\n// Paragraph 3: duplicate \"has_unique_chars\" listing -> removed.\n// Paragraph 4: trailing blank paragraph -> left alone.\nconst targetParagraph = paragraphs.items[0];\nconst paragraphsToRemove = [\n  paragraphs.items[3],\n  paragraphs.items[2],\n  paragraphs.items[1],\n];\n\nfor (const p of paragraphsToRemove) {\n  p.delete();\n}\nawait context.sync();\n\n// Clear the first paragraph's content, then inject the new run structure\n// (one run per source line, each run carrying a leading <w:br/> plus the\n// line's text) via raw OOXML so the resulting markup mirrors a real Word\n// edit rather than collapsing into a single merged run.\ntargetParagraph.clear();\nawait context.sync();\n\nconst newParagraphOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:br/><w:t>class bank():</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">    def bank_details(self):</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        print(\"bank name = SBI Ernakulam\")</w:t></w:r>\n            <w:r><w:br/><w:t>class client(bank):</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">    def __init__(self,name,account_no,balance):</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        self.name = name</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        self.account_number = account_no</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        self.balance = balance</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">    def b_details(self):</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        print(self.name, \",name:\")</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        print(self.account_number, \",balance:\")</w:t></w:r>\n            <w:r><w:br/><w:t xml:space=\"preserve\">        print(self.account_number, \",account number:\")</w:t></w:r>\n            <w:r><w:br/></w:r>\n            <w:r><w:br/></w:r>\n            <w:r><w:br/></w:r>\n            <w:r><w:br/><w:t>B1=client(\"akhil\", \"0098765\", \"1000rs\")</w:t></w:r>\n            <w:r><w:br/><w:t>B2=client(\"anu\", \"00954321\", \"500rs\")</w:t></w:r>\n            <w:r><w:br/><w:t>#</w:t></w:r>\n            <w:r><w:br/><w:t>B1.b_details()</w:t></w:r>\n            <w:r><w:br/><w:t>B2.b_details()</w:t></w:r>\n            <w:r><w:br/></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst insertedRange = targetParagraph.getRange(\"Start\");\ninsertedRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the first code listing (\"has_unique_chars\") with the new\n# \"bank\"/\"client\" class listing, and drop the two blank paragraphs plus\n# the duplicate \"has_unique_chars\" listing that followed it, leaving the\n# trailing blank paragraph untouched.\n\n$d = $word.ActiveDocument\n\nWrite-Output (\"Paragraph count before: \" + $d.Paragraphs.Count)\n\n# Paragraph 1: first \"has_unique_chars\" listing -> becomes the new bank code.\n# Paragraph 2: blank paragraph -> removed.\n# Paragraph 3: blank paragraph -> removed.\n# Paragraph 4: duplicate \"has_unique_chars\" listing -> removed.\n# Paragraph 5: trailing blank paragraph -> left alone.\n# Delete from the back so earlier paragraph indices stay valid.\n$d.Paragraphs(4).Range.Delete() | Out-Null\n$d.Paragraphs(3).Range.Delete() | Out-Null\n$d.Paragraphs(2).Range.Delete() | Out-Null\n\nWrite-Output (\"Paragraph count after removal: \" + $d.Paragraphs.Count)\n\n# Replace paragraph 1's content (minus its trailing paragraph mark) with the\n# new bank/client listing, injected as raw WordprocessingML so each source\n# line becomes its own run (leading <w:br/> + the line text), matching how\n# Word itself would represent a pasted multi-line listing.\n$target = $d.Paragraphs(1).Range\n$target.MoveEnd(1, -1) | Out-Null\n$target.Delete() | Out-Null\n\n$newRunsXml = '<w:r><w:br/><w:t>class bank():</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">    def bank_details(self):</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        print(\"bank name = SBI Ernakulam\")</w:t></w:r><w:r><w:br/><w:t>class client(bank):</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">    def __init__(self,name,account_no,balance):</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        self.name = name</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        self.account_number = account_no</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        self.balance = balance</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">    def b_details(self):</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        print(self.name, \",name:\")</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        print(self.account_number, \",balance:\")</w:t></w:r><w:r><w:br/><w:t xml:space=\"preserve\">        print(self.account_number, \",account number:\")</w:t></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/></w:r><w:r><w:br/><w:t>B1=client(\"akhil\", \"0098765\", \"1000rs\")</w:t></w:r><w:r><w:br/><w:t>B2=client(\"anu\", \"00954321\", \"500rs\")</w:t></w:r><w:r><w:br/><w:t>#</w:t></w:r><w:r><w:br/><w:t>B1.b_details()</w:t></w:r><w:r><w:br/><w:t>B2.b_details()</w:t></w:r><w:r><w:br/></w:r>'\n\n$newParagraphXml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' + $newRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$d.Paragraphs(1).Range.InsertXML($newParagraphXml) | Out-Null\n\nWrite-Output (\"Paragraph count final: \" + $d.Paragraphs.Count)\nWrite-Output (\"Paragraph 1 text: \" + $d.Paragraphs(1).Range.Text)\n"}
